$wb = $excel.ActiveWorkbook

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16
$ws.Range("H5").Value = 291.42856
$ws.Range("I5").Value = 276
$ws.Range("J5").Value = 330
$ws.Range("K5").Value = 276
$ws.Range("L5").Value = 330
$ws.Range("M5").Value = -164
$ws.Range("N5").Value = -554
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H37").Value = 20609.945
$ws.Range("I37").Value = 7500
$ws.Range("J37").Value = 22248.688
$ws.Range("K37").Value = 7500
$ws.Range("L37").Value = 22248.688
$ws.Range("M37").Value = -7227
$ws.Range("N37").Value = -22794.688
$ws.Range("H44").Value = 27169
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 27169
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 27169
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -28145
$ws.Range("H55").Value = 16466
$ws.Range("I55").Value = 10047
$ws.Range("J55").Value = 22885
$ws.Range("K55").Value = 10047
$ws.Range("L55").Value = 22885
$ws.Range("M55").Value = -9732
$ws.Range("N55").Value = -23515
$ws.Range("H63").Value = 2927.9697
$ws.Range("I63").Value = 2606.8333
$ws.Range("J63").Value = 3313.3333
$ws.Range("K63").Value = 2606.8333
$ws.Range("L63").Value = 3313.3333
$ws.Range("M63").Value = -1920.8333
$ws.Range("N63").Value = -4685.3333
$ws.Range("H66").Value = 2927.9697
$ws.Range("I66").Value = 2606.8333
$ws.Range("J66").Value = 3313.3333
$ws.Range("K66").Value = 13034.1665
$ws.Range("L66").Value = 16566.6665
$ws.Range("M66").Value = -9602.166499999999
$ws.Range("N66").Value = -23430.6665
$ws.Range("H80").Value = 30061.111
$ws.Range("J80").Value = 30061.111
$ws.Range("L80").Value = 30061.111
$ws.Range("N80").Value = -32057.111
$ws.Range("H83").Value = 30061.111
$ws.Range("J83").Value = 30061.111
$ws.Range("L83").Value = 90183.333
$ws.Range("N83").Value = -100167.333
$ws.Range("H102").Value = 8442.637000000001
$ws.Range("I102").Value = 5306.9
$ws.Range("J102").Value = 39800
$ws.Range("K102").Value = 5306.9
$ws.Range("L102").Value = 39800
$ws.Range("M102").Value = -3684.9
$ws.Range("N102").Value = -43044

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 291.42856
$ws.Range("I4").Value = 276
$ws.Range("J4").Value = 330
$ws.Range("K4").Value = 276
$ws.Range("L4").Value = 330
$ws.Range("M4").Value = -161
$ws.Range("N4").Value = -560
$ws.Range("H15").Value = 22862.143
$ws.Range("J15").Value = 22862.143
$ws.Range("L15").Value = 22862.143
$ws.Range("N15").Value = -23316.143
$ws.Range("H19").Value = 18341.666
$ws.Range("J19").Value = 18341.666
$ws.Range("L19").Value = 18341.666
$ws.Range("N19").Value = -18687.666
$ws.Range("H35").Value = 21874
$ws.Range("J35").Value = 21874
$ws.Range("L35").Value = 21874
$ws.Range("N35").Value = -22494
$ws.Range("H82").Value = 14839.7
$ws.Range("I82").Value = 1330.5714
$ws.Range("J82").Value = 46361
$ws.Range("K82").Value = 1330.5714
$ws.Range("L82").Value = 46361
$ws.Range("M82").Value = -947.5714
$ws.Range("N82").Value = -47127
$ws.Range("H85").Value = 14839.7
$ws.Range("I85").Value = 1330.5714
$ws.Range("J85").Value = 46361
$ws.Range("K85").Value = 1330.5714
$ws.Range("L85").Value = 46361
$ws.Range("M85").Value = -4.57140000000004
$ws.Range("N85").Value = -49013

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5246.2
$ws.Range("I31").Value = 6579.5557
$ws.Range("J31").Value = 4496.1875
$ws.Range("K31").Value = 6579.5557
$ws.Range("L31").Value = 4496.1875
$ws.Range("M31").Value = -6284.5557
$ws.Range("N31").Value = -5086.1875
$ws.Range("H34").Value = 5246.2
$ws.Range("I34").Value = 6579.5557
$ws.Range("J34").Value = 4496.1875
$ws.Range("K34").Value = 6579.5557
$ws.Range("L34").Value = 4496.1875
$ws.Range("M34").Value = -6377.5557
$ws.Range("N34").Value = -4900.1875
$ws.Range("H41").Value = 12198.667
$ws.Range("J41").Value = 19720.834
$ws.Range("L41").Value = 19720.834
$ws.Range("N41").Value = -20576.834
$ws.Range("H50").Value = 20390.834
$ws.Range("J50").Value = 20390.834
$ws.Range("L50").Value = 20390.834
$ws.Range("N50").Value = -21640.834
$ws.Range("H51").Value = 19138.416
$ws.Range("J51").Value = 19514.637
$ws.Range("L51").Value = 19514.637
$ws.Range("N51").Value = -20986.637
$ws.Range("H59").Value = 34282.832
$ws.Range("J59").Value = 34282.832
$ws.Range("L59").Value = 34282.832
$ws.Range("N59").Value = -36572.832
$ws.Range("H60").Value = 19142.666
$ws.Range("I60").Value = 15950
$ws.Range("J60").Value = 19781.2
$ws.Range("K60").Value = 15950
$ws.Range("L60").Value = 19781.2
$ws.Range("M60").Value = -15439
$ws.Range("N60").Value = -20803.2
$ws.Range("H61").Value = 19138.416
$ws.Range("J61").Value = 19514.637
$ws.Range("L61").Value = 19514.637
$ws.Range("N61").Value = -20210.637
$ws.Range("H68").Value = 28755.65
$ws.Range("I68").Value = 9000
$ws.Range("J68").Value = 29795.422
$ws.Range("K68").Value = 9000
$ws.Range("L68").Value = 29795.422
$ws.Range("M68").Value = -8251
$ws.Range("N68").Value = -31293.422
$ws.Range("H71").Value = 28755.65
$ws.Range("I71").Value = 9000
$ws.Range("J71").Value = 29795.422
$ws.Range("K71").Value = 27000
$ws.Range("L71").Value = 89386.266
$ws.Range("M71").Value = -23256
$ws.Range("N71").Value = -96874.266
$ws.Range("H74").Value = 29680.215
$ws.Range("J74").Value = 29680.215
$ws.Range("L74").Value = 29680.215
$ws.Range("N74").Value = -31428.215
$ws.Range("H77").Value = 29680.215
$ws.Range("J77").Value = 29680.215
$ws.Range("L77").Value = 89040.645
$ws.Range("N77").Value = -97776.645

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1717.4615
$ws.Range("I113").Value = 2530.7144
$ws.Range("J113").Value = 768.6667
$ws.Range("K113").Value = 7592.1432
$ws.Range("L113").Value = 2306.0001
$ws.Range("M113").Value = -5422.1432
$ws.Range("N113").Value = -6646.0001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1679.3334
$ws.Range("I132").Value = 1279.1052
$ws.Range("K132").Value = 3837.3156
$ws.Range("M132").Value = -1307.3156

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4490.4575
$ws.Range("I132").Value = 3020.8572
$ws.Range("K132").Value = 9062.571599999999
$ws.Range("M132").Value = -6532.571599999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 197961.1
$ws.Range("I136").Value = 226012.22
$ws.Range("J136").Value = 1603.25
$ws.Range("K136").Value = 678036.66
$ws.Range("L136").Value = 4809.75
$ws.Range("M136").Value = -675486.66
$ws.Range("N136").Value = -9909.75
